# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 3.082599426703578;  E = 0.4998867070740569; G = 8.418600821238126 }
    3 = @{ B = 0.7287194209349384; C = 0.05231270169004087; D = 0.1529057820181812; E = 0.4998867070740569; G = 1.433824611717217 }
    4 = @{ B = 0.3464964993005633; C = 0.3375848360084654;  D = 3.082599426703578;  E = 0.4998867070740569; G = 4.266567469086664 }
    5 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 3.082599426703578;  E = 0.4998867070740569; G = 8.418600821238126 }
    6 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
    7 = @{ B = 0.06328177979961902; C = 0.3375848360084654; D = 0.7127328510149897; E = 0.4998867070740569; G = 1.613486173897131 }
    8 = @{ B = 0.1554434735375247; C = 0.05231270169004087; D = 0.7127328510149897; E = 0.4998867070740569; G = 1.420375733316612 }
}

foreach ($row in $data.Keys) {
    $rowVals = $data[$row]
    $ws.Range("B$row").Value = $rowVals.B
    $ws.Range("C$row").Value = $rowVals.C
    $ws.Range("D$row").Value = $rowVals.D
    $ws.Range("E$row").Value = $rowVals.E
    $ws.Range("G$row").Value = $rowVals.G
}
